$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login_sheet")

# Update password test value in row 4 (was "Test123", now "Test")
$ws.Range("B4").Value = "Test"

# Update the active selection to reflect the latest cursor position
$ws.Range("L10").Select()
